# "Fruta / hortaliza, semanal" — insert a new weekly price-report row for
# Papa (Red Lady, Región del Bíobío) ahead of the existing row 131, pushing
# every following record down by one row (131-189 -> 132-190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 131; Excel shifts rows 131..189 down
# to 132..190 and the sheet's used range/dimension grow to R190 automatically.
$ws.Rows.Item(131).Insert()

# Populate the newly-inserted row 131 with the new record.
$ws.Range("A131").Value2 = 1
$ws.Range("B131").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C131").Value2 = "Arica y Parinacota"
$ws.Range("D131").Value2 = 45027
$ws.Range("E131").Value2 = 15
$ws.Range("F131").Value2 = 100114001
$ws.Range("G131").Value2 = "Papa"
$ws.Range("H131").Value2 = "Red Lady"
$ws.Range("I131").Value2 = "1a (cosecha)"
$ws.Range("J131").Value2 = 1100
$ws.Range("K131").Value2 = 13000
$ws.Range("L131").Value2 = 14000
$ws.Range("M131").Value2 = 13409
$ws.Range("N131").Value2 = "$/saco 25 kilos"
$ws.Range("O131").Value2 = "Región del Bíobío"
$ws.Range("P131").Value2 = 536
$ws.Range("Q131").Value2 = 25
$ws.Range("R131").Value2 = "Hortaliza"

# Match the date cell style (yyyy-mm-dd style number format) used by every
# other row's date column.
$ws.Range("D131").NumberFormat = $ws.Range("D132").NumberFormat
